# Insert a new weekly price record as the first data row (row 488),
# pushing the existing rows 488:588 down to 489:589.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(488).Insert()

$ws.Cells.Item(488, 1).Value = 3
$ws.Cells.Item(488, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(488, 3).Value = "Coquimbo"
$ws.Cells.Item(488, 4).Value = 44889
$ws.Cells.Item(488, 5).Value = 5
$ws.Cells.Item(488, 6).Value = 100112037
$ws.Cells.Item(488, 7).Value = "Cebollín"
$ws.Cells.Item(488, 8).Value = "Sin especificar"
$ws.Cells.Item(488, 9).Value = "Primera"
$ws.Cells.Item(488, 10).Value = 250
$ws.Cells.Item(488, 11).Value = 3500
$ws.Cells.Item(488, 12).Value = 3800
$ws.Cells.Item(488, 13).Value = 3644
$ws.Cells.Item(488, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(488, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(488, 16).Value = 101
$ws.Cells.Item(488, 17).Value = 36
$ws.Cells.Item(488, 18).Value = "Hortaliza"
